$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Analisis y concepción" ---
$ws1 = $wb.Worksheets.Item("Analisis y concepción")
$ws1.Range("B2").Value = 1
$ws1.Range("B3").Value = 1
$ws1.Range("B4").Value = 1
$ws1.Range("B5").Value = 1
$ws1.Range("B6").Value = 1
$ws1.Range("B7").Select()

# --- Sheet 2: "Desarrollo de la aplicación" ---
$ws2 = $wb.Worksheets.Item("Desarrollo de la aplicación")
$ws2.Range("C2").Value = 1
$ws2.Range("C3").Value = 1
$ws2.Range("C5").Value = 1

$ws2.Range("B7").Value = "front-end slider"
$ws2.Range("C7").Value = 1
$ws2.Range("B8").Value = "front galeria de fondos miembros"
$ws2.Range("B9").Value = "front galeria de fondos publico"
$ws2.Range("B10").Value = "front end vista de galeria con overlay"
$ws2.Range("B11").Value = "back end crud de menus hubo cambios"

$ws2.Activate()
$ws2.Range("C12").Select()
